$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.72 = 14256.51 pesos`n✅ 14256.51 pesos = 3.7 = 937.34 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the N10/O10/N12/O12 rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 269
$wsTasas.Range("O10").Value = 3835
$wsTasas.Range("N12").Value = 3850
$wsTasas.Range("O12").Value = 253.13
